$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "最大销售利润"
$ws.Range("J1").Value = 98790
$ws.Range("I2").Value = "最小销售利润"
$ws.Range("J2").Value = 27750

$ws.Range("J1:J2").NumberFormat = '"¥"#,##0.00;"¥"\-#,##0.00'

$ws.Columns.AutoFit() | Out-Null
